# Apply the "uploader with image path" update to row 2 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# category_ids / sku
$ws.Range("A2").Value = 19
$ws.Range("B2").Value = 4554

# image / small_image / thumbnail now point at the full uploader URL
$newImage = "https://www.pdcorders.com/media/import/pdc.png"
$ws.Range("K2").Value = $newImage
$ws.Range("L2").Value = $newImage
$ws.Range("M2").Value = $newImage

# name / product_name / description / short_description all share the
# same updated product title text
$newName = "ELEGANCE GOLD LONG LASTING LIP STICK 9 MAROON"
$ws.Range("C2").Value = $newName
$ws.Range("D2").Value = $newName
$ws.Range("I2").Value = $newName
$ws.Range("J2").Value = $newName

# status / visibility become text values instead of numeric codes
$ws.Range("U2").Value = "Enabled"
$ws.Range("V2").Value = "Catalog, Search"
